# Elimina dato 2023 de población censada
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 7 (Censos, 2023, 3444.263); remaining rows shift up
$ws.Rows.Item(7).Delete()

# Refresh the stale sort-state range (was A2:C14 / A2:A14 / B2:B14) now that
# the data only spans through row 13
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A2:A13"))
$sort.SortFields.Add($ws.Range("B2:B13"))
$sort.SetRange($ws.Range("A2:C13"))
$sort.Header = 0
$sort.Apply()

# Refresh the stale _FilterDatabase defined name (was $A$1:$C$14)
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Hoja1!`$A`$1:`$C`$13"
    }
}

# Move selection to reflect the saved cursor position in the authored file
$ws.Range("D17").Select()
